# Insert three brand-new price rows above the current row 905 of the
# "Vega Modelo de Temuco - Uva" data table. This pushes all existing rows
# 905-970 down to 908-973 (dimension grows from A1:T970 to A1:T973) and
# Excel copies the row-above formatting (incl. the date style on column D).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("905:907").Insert()

$newRows = @(
    @{ Row = 905; Values = @(10, "Vega Modelo de Temuco", "La Araucanía", 44769, 9, "Fruta", 100109, "Uva", 100109001, "Uva", "Autumn Royal", "Primera", 250, 6000, 6000, 6000, "`$/bandeja 8 kilos", "Región de O'Higgins", 750, 8) },
    @{ Row = 906; Values = @(10, "Vega Modelo de Temuco", "La Araucanía", 44769, 9, "Fruta", 100109, "Uva", 100109001, "Uva", "Crimpson Seedless", "Primera", 380, 8000, 8000, 8000, "`$/bandeja 8 kilos", "Región de O'Higgins", 1000, 8) },
    @{ Row = 907; Values = @(10, "Vega Modelo de Temuco", "La Araucanía", 44340, 9, "Fruta", 100109, "Uva", 100109001, "Uva", "Red Globe", "Primera", 255, 12000, 12000, 12000, "`$/bandeja 10 kilos", "Región de O'Higgins", 1200, 10) }
)

foreach ($entry in $newRows) {
    $r = $entry.Row
    $vals = $entry.Values
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($r, $i + 1).Value = $vals[$i]
    }
}
